$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: header labels (bold/centered style already applied via s="1")
$ws.Range("A1").Value = "Ders Çıktı"
$ws.Range("B1").Value = "TABLO 3"
$ws.Range("C1").Value = "Ağırlıklı Değerlendirme"
$ws.Range("D1").Value = "Ağırlıklı Değerlendirme"
$ws.Range("E1").Value = "Ağırlıklı Değerlendirme"
$ws.Range("F1").Value = "Ağırlıklı Değerlendirme"
$ws.Range("G1").Value = "Toplam"

# Row 2: column headers
$ws.Range("A2").Value = "Ders Çıktı"
$ws.Range("B2").Value = "Öd1"
$ws.Range("C2").Value = "Öd2"
$ws.Range("D2").Value = "Quiz"
$ws.Range("E2").Value = "Vize"
$ws.Range("F2").Value = "Fin"
$ws.Range("G2").Value = "Toplam"

# Row 3
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 0.1
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0.1
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0.2

# Row 4
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 0.1
$ws.Range("C4").Value = 0.1
$ws.Range("D4").Value = 0.1
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0.4
$ws.Range("G4").Value = 0.7000000000000001

# Row 5
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0.1
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0.3
$ws.Range("F5").Value = 0.4
$ws.Range("G5").Value = 0.8

# Row 6
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0.3
$ws.Range("F6").Value = 0.4
$ws.Range("G6").Value = 0.7

# Row 7
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0.1
$ws.Range("E7").Value = 0.3
$ws.Range("F7").Value = 0.4
$ws.Range("G7").Value = 0.8

# Remove the now-unused row 8 (table shrank from 5 data rows + headers = 8 rows to 7 rows)
$ws.Rows.Item(8).Delete()
